$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "таблицой" -> "таблицей" (C5, C14)
$ws.Range("C5").Value = "1. Настроить программу lab-grader на работу с тестовой таблицей: https://docs.google.com/spreadsheets/d/17Qy4CFnqB3jKV9WC1TYzxQSlK-tolBmrW4Y8m1w2QXg/edit#gid=1593948738"
$ws.Range("C14").Value = "1. Настроить программу lab-grader на работу с тестовой таблицей: https://docs.google.com/spreadsheets/d/17Qy4CFnqB3jKV9WC1TYzxQSlK-tolBmrW4Y8m1w2QXg/edit#gid=1593948738"

# Fix typo: "пунке" -> "пункте" (C7 only; C16 gets a different typo "пунте")
$ws.Range("C7").Value = "2. После того как в пункте со статусом выполнения тестов на сайте: https://github.com/suai-ms-2020/ms-task1-BatMaxim появиться галочка, в консоли выполнения теста нажать любую клавишу для продолжения работы теста."
$ws.Range("C16").Value = "2. После того как в пунте со статусом выполнения тестов на сайте: https://github.com/suai-ms-2020/ms-task1-BatMaxim появиться галочка, в консоли выполнения теста нажать любую клавишу для продолжения работы теста."

# Fix typo: "репозиторя" -> "репозитория" (D8, D17)
$ws.Range("D8").Value = "Из удаленного репозитория студента будут удалены файлы проекта."
$ws.Range("D17").Value = "Из удаленного репозитория студента будут удалены файлы проекта."

# Fix typos in D9 ("ячейче"->"ячейке", "должа"->"должна")
$ws.Range("D9").Value = "В ячейке соответствующей 1й лабораторной работе для студента Иванова Ивана Ивановича должна появиться галочка."

# D18: fix typos AND change text to describe the invalid-case expected result
$ws.Range("D18").Value = "В ячейке соответствующей 1й лабораторной работе для студента Иванова Ивана Ивановича должна появиться надпись, сообщающая о неправильном варианте."

# Row 18 grew taller because its wrapped text got longer
$ws.Range("B18").RowHeight = 90

# Update sheet view (scroll position / selection) to match the saved state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H14").Select() | Out-Null
